$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B29's formula to add an extra hour (6+2+2 -> 6+2+2+1)
$ws.Range("B29").Formula = "=6+2+2+1"

# Force recalculation so dependent SUM/AVERAGE formulas update
$excel.Calculate()

# Move the active selection to B30 (matching the saved selection in the file)
$ws.Range("B30").Select()
